# RDCC-5182 Added Version check
# Adds a new "VERSION" worksheet containing a small "File version" / "vx.xx"
# table (placed at row 6, columns A:B), positions it as the second tab, and
# makes it the active sheet.

$wb = $excel.ActiveWorkbook

# Worksheets.Add() inserts the new sheet before the currently active sheet
# and makes it active, so populate it with its data first (while it is the
# active sheet) and only move/re-activate afterwards.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "VERSION"

# Populate the version info, starting at row 6 as in the target workbook.
$newSheet.Range("A6").Value = "File version"
$newSheet.Range("B6").Value = "vx.xx"
$newSheet.Range("B6").Select()

# Move the new sheet after the existing one so it becomes the second tab.
$existing = $wb.Worksheets.Item("Service to CW Roles Mapping")
$newSheet.Move($null, $existing)

# Make the VERSION sheet the active/visible tab.
$wb.Worksheets.Item("VERSION").Activate()
